$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 5 is the "donor_PHM_calc" row. Update the summary statistics from
# median [IQR] to mean ± SD, update the p-value, and change the test
# from Wilcoxon rank-sum to t-test.
$t.Cell(5, 2).Range.Text = "190.6 ± 33.6"
$t.Cell(5, 3).Range.Text = "196.7 ± 37.5"
$t.Cell(5, 4).Range.Text = "191.4 ± 33.8"
$t.Cell(5, 5).Range.Text = "0.676"
$t.Cell(5, 6).Range.Text = "t-test"
